# Adds new 3D-print request log entries for "CAS" (added by Josh on 24-05-2018),
# matching commit: "added new entries for prints for CAS".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Date Requested | B Date Completed | C Requestor | D Part |
#          E Quantity | F Material | G Shell Count (#) | H Infill (%) |
#          I Layer Height (mm) | J Comments
# (Note: column B - Date Completed - is intentionally left blank for these rows,
# matching the source data.)
$newEntries = @(
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS a4c 100 scale";        Qty = 2;  Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS a4c 50 scale";         Qty = 20; Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS LAX 100 scale";        Qty = 2;  Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS LAX 50 scale";         Qty = 20; Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS SAX 100 scale";        Qty = 2;  Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "FOCUS SAX 50 scale";         Qty = 20; Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" },
    @{ Date = "24-05-2018"; Requestor = "Josh"; Part = "Modular Bloodpool Model";    Qty = 10; Material = "Polylite"; Shells = 2; Infill = 20; Layer = 0.2; Comments = "For CAS" }
)

$startRow = 20
$row = $startRow
foreach ($entry in $newEntries) {
    $ws.Range("A$row").Value = $entry.Date
    $ws.Range("C$row").Value = $entry.Requestor
    $ws.Range("D$row").Value = $entry.Part
    $ws.Range("E$row").Value = $entry.Qty
    $ws.Range("F$row").Value = $entry.Material
    $ws.Range("G$row").Value = $entry.Shells
    $ws.Range("H$row").Value = $entry.Infill
    $ws.Range("I$row").Value = $entry.Layer
    $ws.Range("J$row").Value = $entry.Comments
    $row++
}

# A couple of stray rows below the new data were touched (e.g. by scrolling /
# clicking around) while entering the new entries, leaving them present with
# their (auto) row height but no cell content.
$ws.Range("A41").EntireRow.RowHeight = 13.8
$ws.Range("A43").EntireRow.RowHeight = 13.8

# Scroll the frozen view down so row 14 is the first visible row below the
# frozen header, and leave the active cell where the user last left off.
$av = $excel.ActiveWindow
$av.ScrollRow = 14
[void]$ws.Range("D28").Select()
